# Apply Anders Meidahl's reviewed edits:
#  - accept the two tracked "noget " deletions in the first paragraph
#  - move the auto "_GoBack" bookmark from its old cursor position
#    (between "mest værdifulde. " and "For en gennemgang...") to the
#    new last-edit position (between "hardware og " and "software.")
#  - let Word's run normalisation merge the runs that become adjacent
#    once the deleted text / bookmark no longer separates them

$d = $word.ActiveDocument

$origTrack = $d.TrackRevisions
$d.TrackRevisions = $false

# Accept every still-pending tracked change (the two "noget " deletions).
# Accepted one at a time (rather than AcceptAll) so Word doesn't also
# normalise/strip the unrelated rsid bookkeeping on runs elsewhere in
# the document.
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions($i).Accept()
}

# Relocate the "_GoBack" bookmark: delete the stale one (old cursor spot)
# and drop a fresh, empty one where the last accepted edit happened.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$anchor = $d.Content
$anchor.Find.Execute("hardware og ", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$goBackRange = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Nudge Word into normalising (merging) the runs left/right of each
# removed deletion so the text collapses back into single <w:r> runs,
# matching what Word itself does when the markup is saved after review.
$m1 = $d.Content
$m1.Find.Execute("Fælles", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Fælles", 2)

$m2 = $d.Content
$m2.Find.Execute("mest værdifulde", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "mest værdifulde", 2)

$d.TrackRevisions = $origTrack
